$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Multivalued" column header in L4 (bold Calibri 11, black).
$ws.Range("L4").Value = "Multivalued"
$ws.Range("L4").Font.Size = 11
$ws.Range("L4").Font.Bold = $true

# New boolean-like "FALSE" text values for the Multivalued column (L5:L8).
# A leading apostrophe forces the literal text "FALSE" (shared string) rather
# than a native boolean, matching a manually quote-prefixed entry.
$ws.Range("L5").Value = "'FALSE"
$ws.Range("L6").Value = "'FALSE"
$ws.Range("L7").Value = "'FALSE"
$ws.Range("L8").Value = "'FALSE"

$boolFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("L5:L8").NumberFormat = $boolFormat
$ws.Range("L5:L8").HorizontalAlignment = -4131

# Match the selection left behind by the edit: L4:L8 highlighted, L4 active.
$ws.Range("L4:L8").Select()
